$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4931.048
$ws.Range("I137").Value = 3070.5173
$ws.Range("K137").Value = 9211.5519
$ws.Range("M137").Value = -6661.5519

$ws.Range("H138").Value = 2942.9841
$ws.Range("I138").Value = 1924.7333
$ws.Range("J138").Value = 3261.1875
$ws.Range("K138").Value = 5774.199900000001
$ws.Range("L138").Value = 9783.5625
$ws.Range("M138").Value = -634.1999000000005
$ws.Range("N138").Value = -20063.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null

$ws.Range("H39").Value = 15585.5
$ws.Range("J39").Value = 31999
$ws.Range("L39").Value = 31999
$ws.Range("N39").Value = -33039

$ws.Range("H41").Value = 5554.8887
$ws.Range("I41").Value = 5554.8887
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5554.8887
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -5140.8887
$ws.Range("N41").Value = $null

$ws.Range("H42").Value = 2001
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2001
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 2001
$ws.Range("M42").Value = $null
$ws.Range("N42").Value = -2973

$ws.Range("H61").Value = 26793138
$ws.Range("I61").Value = 25006524
$ws.Range("J61").Value = 31259672
$ws.Range("K61").Value = 25006524
$ws.Range("L61").Value = 31259672
$ws.Range("M61").Value = -25006312
$ws.Range("N61").Value = -31260096

$ws.Range("H74").Value = 5656789.5
$ws.Range("I74").Value = 8067110.5
$ws.Range("J74").Value = 675459.4
$ws.Range("K74").Value = 8067110.5
$ws.Range("L74").Value = 675459.4
$ws.Range("M74").Value = -8066236.5
$ws.Range("N74").Value = -677207.4

$ws.Range("H77").Value = 5656789.5
$ws.Range("I77").Value = 8067110.5
$ws.Range("J77").Value = 675459.4
$ws.Range("K77").Value = 40335552.5
$ws.Range("L77").Value = 3377297
$ws.Range("M77").Value = -40331184.5
$ws.Range("N77").Value = -3386033

$ws.Range("H132").Value = 5187.727
$ws.Range("I132").Value = 2012
$ws.Range("K132").Value = 6036
$ws.Range("M132").Value = -3506

$ws.Range("H136").Value = 26793138
$ws.Range("I136").Value = 25006524
$ws.Range("J136").Value = 31259672
$ws.Range("K136").Value = 75019572
$ws.Range("L136").Value = 93779016
$ws.Range("M136").Value = -75017022
$ws.Range("N136").Value = -93784116

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4556.1514
$ws.Range("I20").Value = 5458.15
$ws.Range("J20").Value = 3168.4614
$ws.Range("K20").Value = 5458.15
$ws.Range("L20").Value = 3168.4614
$ws.Range("M20").Value = -5211.15
$ws.Range("N20").Value = -3662.4614

$ws.Range("H99").Value = 8052.2856
$ws.Range("I99").Value = 14821.625
$ws.Range("K99").Value = 14821.625
$ws.Range("M99").Value = -13323.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13480427
$ws.Range("I31").Value = 34496396
$ws.Range("K31").Value = 34496396
$ws.Range("M31").Value = -34496101

$ws.Range("H34").Value = 13480427
$ws.Range("I34").Value = 34496396
$ws.Range("K34").Value = 34496396
$ws.Range("M34").Value = -34496194

$ws.Range("H38").Value = 1521
$ws.Range("J38").Value = 1521
$ws.Range("L38").Value = 1521
$ws.Range("N38").Value = -2275

$ws.Range("H46").Value = 1521
$ws.Range("J46").Value = 1521
$ws.Range("L46").Value = 1521
$ws.Range("N46").Value = -1943

$ws.Range("H94").Value = 3496.9546
$ws.Range("J94").Value = 3313.4119
$ws.Range("L94").Value = 3313.4119
$ws.Range("N94").Value = -4215.4119

$ws.Range("H132").Value = 3079.4546
$ws.Range("I132").Value = 2730.8333
$ws.Range("K132").Value = 8192.499899999999
$ws.Range("M132").Value = -5662.499899999999

$ws.Range("H134").Value = 1695.7142
$ws.Range("I134").Value = 1233.1945
$ws.Range("K134").Value = 3699.5835
$ws.Range("M134").Value = -1164.5835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 628.6667
$ws.Range("I11").Value = 692.5
$ws.Range("J11").Value = 501
$ws.Range("K11").Value = 2077.5
$ws.Range("L11").Value = 1503
$ws.Range("M11").Value = -1937.5
$ws.Range("N11").Value = -1783

$ws.Range("H20").Value = 8899
$ws.Range("J20").Value = 8899
$ws.Range("L20").Value = 26697
$ws.Range("N20").Value = -27151

$ws.Range("H23").Value = 1385.7142
$ws.Range("I23").Value = 1385.7142
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4157.142599999999
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3922.142599999999
$ws.Range("N23").Value = $null

$ws.Range("H26").Value = 130.33333
$ws.Range("I26").Value = 106.625
$ws.Range("K26").Value = 319.875
$ws.Range("M26").Value = -31.875

$ws.Range("H33").Value = 307
$ws.Range("I33").Value = 287.25
$ws.Range("J33").Value = 333.33334
$ws.Range("K33").Value = 1723.5
$ws.Range("L33").Value = 2000.00004
$ws.Range("M33").Value = -1440.5
$ws.Range("N33").Value = -2566.00004

$ws.Range("I35").Value = 300
$ws.Range("K35").Value = 900
$ws.Range("M35").Value = -612

$ws.Range("H38").Value = 184.78572
$ws.Range("I38").Value = 508.75
$ws.Range("J38").Value = 55.2
$ws.Range("K38").Value = 1526.25
$ws.Range("L38").Value = 165.6
$ws.Range("M38").Value = -1179.25
$ws.Range("N38").Value = -859.6

$ws.Range("H46").Value = 2268.9
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 2997
$ws.Range("M46").Value = -2906

$ws.Range("H113").Value = 1516.3077
$ws.Range("J113").Value = 1516.3077
$ws.Range("L113").Value = 4548.9231
$ws.Range("N113").Value = -8888.9231

$ws.Range("H131").Value = 4847.1313
$ws.Range("I131").Value = 4857.125
$ws.Range("J131").Value = 4843.5776
$ws.Range("K131").Value = 14571.375
$ws.Range("L131").Value = 14530.7328
$ws.Range("M131").Value = -9531.375
$ws.Range("N131").Value = -24610.7328

$ws.Range("H139").Value = 2314.4092
$ws.Range("I139").Value = 1853.1562
$ws.Range("J139").Value = 3544.4167
$ws.Range("K139").Value = 5559.4686
$ws.Range("L139").Value = 10633.2501
$ws.Range("M139").Value = -419.4686000000002
$ws.Range("N139").Value = -20913.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1860
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 2000
$ws.Range("N31").Value = -2584

$ws.Range("H37").Value = 1860
$ws.Range("J37").Value = 2000
$ws.Range("L37").Value = 2000
$ws.Range("N37").Value = -2554

$ws.Range("H107").Value = 420.25
$ws.Range("I107").Value = 314.875
$ws.Range("J107").Value = 525.625
$ws.Range("K107").Value = 314.875
$ws.Range("L107").Value = 525.625
$ws.Range("M107").Value = 1605.125
$ws.Range("N107").Value = -4365.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58531.105
$ws.Range("J7").Value = 132082.5
$ws.Range("L7").Value = 132082.5
$ws.Range("N7").Value = -132306.5

$ws.Range("H16").Value = 641.0625
$ws.Range("J16").Value = 490.4
$ws.Range("L16").Value = 490.4
$ws.Range("N16").Value = -830.4

$ws.Range("H19").Value = 10375
$ws.Range("J19").Value = 4625
$ws.Range("L19").Value = 4625
$ws.Range("N19").Value = -4965

$ws.Range("H40").Value = 3094.5293
$ws.Range("I40").Value = 2133.9167
$ws.Range("K40").Value = 2133.9167
$ws.Range("M40").Value = -1997.9167

$ws.Range("H55").Value = 37037476
$ws.Range("J55").Value = 530.7778
$ws.Range("L55").Value = 530.7778
$ws.Range("N55").Value = -876.7778

$ws.Range("H75").Value = 99999
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 99999
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 99999
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -101871

$ws.Range("H78").Value = 99999
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 99999
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 299997
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -309357

$ws.Range("H122").Value = 5325.4287
$ws.Range("I122").Value = 4543.1665
$ws.Range("K122").Value = 13629.4995
$ws.Range("M122").Value = -11179.4995

$ws.Range("H126").Value = 58531.105
$ws.Range("J126").Value = 132082.5
$ws.Range("L126").Value = 396247.5
$ws.Range("N126").Value = -401187.5

$ws.Range("H132").Value = 484445.66
$ws.Range("I132").Value = 13168.091
$ws.Range("J132").Value = 1002851
$ws.Range("K132").Value = 39504.273
$ws.Range("L132").Value = 3008553
$ws.Range("M132").Value = -36974.273
$ws.Range("N132").Value = -3013613

$ws.Range("H136").Value = 34496.395
$ws.Range("I136").Value = 4651.241
$ws.Range("K136").Value = 13953.723
$ws.Range("M136").Value = -11403.723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 26373.75
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4851

$ws.Range("H43").Value = 58296.668
$ws.Range("I43").Value = 58296.668
$ws.Range("K43").Value = 58296.668
$ws.Range("M43").Value = -58147.668

$ws.Range("H100").Value = 1153.0869
$ws.Range("I100").Value = 1214.2632
$ws.Range("K100").Value = 2428.5264
$ws.Range("M100").Value = -1887.5264

$ws.Range("H122").Value = 1133.7778
$ws.Range("I122").Value = 1088
$ws.Range("K122").Value = 3264
$ws.Range("M122").Value = -814
